$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("begroting")
$ws1.Activate()
$excel.Goto($ws1.Range("A13"), $true)
$ws1.Range("B40").Select()
